$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: copy style from E1 (bold, border, centered) and set value
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data cells F2:F99: time_taken timestamps (no special style, matches diff)
$ws.Range("F2").Value = "2021-10-05 13:41:51.732895"
$ws.Range("F3").Value = "2021-10-05 13:41:51.732908"
$ws.Range("F4").Value = "2021-10-05 13:41:51.732912"
$ws.Range("F5").Value = "2021-10-05 13:41:51.732915"
$ws.Range("F6").Value = "2021-10-05 13:41:51.732918"
$ws.Range("F7").Value = "2021-10-05 13:41:51.732921"
$ws.Range("F8").Value = "2021-10-05 13:41:51.732924"
$ws.Range("F9").Value = "2021-10-05 13:41:51.732927"
$ws.Range("F10").Value = "2021-10-05 13:41:51.732930"
$ws.Range("F11").Value = "2021-10-05 13:41:51.732933"
$ws.Range("F12").Value = "2021-10-05 13:41:51.732936"
$ws.Range("F13").Value = "2021-10-05 13:41:51.732939"
$ws.Range("F14").Value = "2021-10-05 13:41:51.732942"
$ws.Range("F15").Value = "2021-10-05 13:41:51.732945"
$ws.Range("F16").Value = "2021-10-05 13:41:51.732948"
$ws.Range("F17").Value = "2021-10-05 13:41:51.732951"
$ws.Range("F18").Value = "2021-10-05 13:41:51.732954"
$ws.Range("F19").Value = "2021-10-05 13:41:51.732957"
$ws.Range("F20").Value = "2021-10-05 13:41:51.732960"
$ws.Range("F21").Value = "2021-10-05 13:41:51.732963"
$ws.Range("F22").Value = "2021-10-05 13:41:51.732966"
$ws.Range("F23").Value = "2021-10-05 13:41:51.732968"
$ws.Range("F24").Value = "2021-10-05 13:41:51.732971"
$ws.Range("F25").Value = "2021-10-05 13:41:51.732974"
$ws.Range("F26").Value = "2021-10-05 13:41:51.732978"
$ws.Range("F27").Value = "2021-10-05 13:41:51.732981"
$ws.Range("F28").Value = "2021-10-05 13:41:51.732983"
$ws.Range("F29").Value = "2021-10-05 13:41:51.732986"
$ws.Range("F30").Value = "2021-10-05 13:41:51.732989"
$ws.Range("F31").Value = "2021-10-05 13:41:51.732992"
$ws.Range("F32").Value = "2021-10-05 13:41:51.732995"
$ws.Range("F33").Value = "2021-10-05 13:41:51.732998"
$ws.Range("F34").Value = "2021-10-05 13:41:51.733001"
$ws.Range("F35").Value = "2021-10-05 13:41:51.733004"
$ws.Range("F36").Value = "2021-10-05 13:41:51.733007"
$ws.Range("F37").Value = "2021-10-05 13:41:51.733010"
$ws.Range("F38").Value = "2021-10-05 13:41:51.733013"
$ws.Range("F39").Value = "2021-10-05 13:41:51.733016"
$ws.Range("F40").Value = "2021-10-05 13:41:51.733019"
$ws.Range("F41").Value = "2021-10-05 13:41:51.733022"
$ws.Range("F42").Value = "2021-10-05 13:41:51.733025"
$ws.Range("F43").Value = "2021-10-05 13:41:51.733028"
$ws.Range("F44").Value = "2021-10-05 13:41:51.733031"
$ws.Range("F45").Value = "2021-10-05 13:41:51.733034"
$ws.Range("F46").Value = "2021-10-05 13:41:51.733037"
$ws.Range("F47").Value = "2021-10-05 13:41:51.733040"
$ws.Range("F48").Value = "2021-10-05 13:41:51.733043"
$ws.Range("F49").Value = "2021-10-05 13:41:51.733046"
$ws.Range("F50").Value = "2021-10-05 13:41:51.733049"
$ws.Range("F51").Value = "2021-10-05 13:41:51.733052"
$ws.Range("F52").Value = "2021-10-05 13:41:51.733055"
$ws.Range("F53").Value = "2021-10-05 13:41:51.733058"
$ws.Range("F54").Value = "2021-10-05 13:41:51.733062"
$ws.Range("F55").Value = "2021-10-05 13:41:51.733065"
$ws.Range("F56").Value = "2021-10-05 13:41:51.733067"
$ws.Range("F57").Value = "2021-10-05 13:41:51.733070"
$ws.Range("F58").Value = "2021-10-05 13:41:51.733073"
$ws.Range("F59").Value = "2021-10-05 13:41:51.733076"
$ws.Range("F60").Value = "2021-10-05 13:41:51.733079"
$ws.Range("F61").Value = "2021-10-05 13:41:51.733082"
$ws.Range("F62").Value = "2021-10-05 13:41:51.733085"
$ws.Range("F63").Value = "2021-10-05 13:41:51.733088"
$ws.Range("F64").Value = "2021-10-05 13:41:51.733091"
$ws.Range("F65").Value = "2021-10-05 13:41:51.733094"
$ws.Range("F66").Value = "2021-10-05 13:41:51.733098"
$ws.Range("F67").Value = "2021-10-05 13:41:51.733102"
$ws.Range("F68").Value = "2021-10-05 13:41:51.733104"
$ws.Range("F69").Value = "2021-10-05 13:41:51.733107"
$ws.Range("F70").Value = "2021-10-05 13:41:51.733110"
$ws.Range("F71").Value = "2021-10-05 13:41:51.733113"
$ws.Range("F72").Value = "2021-10-05 13:41:51.733116"
$ws.Range("F73").Value = "2021-10-05 13:41:51.733119"
$ws.Range("F74").Value = "2021-10-05 13:41:51.733122"
$ws.Range("F75").Value = "2021-10-05 13:41:51.733125"
$ws.Range("F76").Value = "2021-10-05 13:41:51.733128"
$ws.Range("F77").Value = "2021-10-05 13:41:51.733131"
$ws.Range("F78").Value = "2021-10-05 13:41:51.733136"
$ws.Range("F79").Value = "2021-10-05 13:41:51.733139"
$ws.Range("F80").Value = "2021-10-05 13:41:51.733142"
$ws.Range("F81").Value = "2021-10-05 13:41:51.733145"
$ws.Range("F82").Value = "2021-10-05 13:41:51.733148"
$ws.Range("F83").Value = "2021-10-05 13:41:51.733151"
$ws.Range("F84").Value = "2021-10-05 13:41:51.733154"
$ws.Range("F85").Value = "2021-10-05 13:41:51.733157"
$ws.Range("F86").Value = "2021-10-05 13:41:51.733160"
$ws.Range("F87").Value = "2021-10-05 13:41:51.733163"
$ws.Range("F88").Value = "2021-10-05 13:41:51.733166"
$ws.Range("F89").Value = "2021-10-05 13:41:51.733169"
$ws.Range("F90").Value = "2021-10-05 13:41:51.733172"
$ws.Range("F91").Value = "2021-10-05 13:41:51.733175"
$ws.Range("F92").Value = "2021-10-05 13:41:51.733178"
$ws.Range("F93").Value = "2021-10-05 13:41:51.733180"
$ws.Range("F94").Value = "2021-10-05 13:41:51.733185"
$ws.Range("F95").Value = "2021-10-05 13:41:51.733188"
$ws.Range("F96").Value = "2021-10-05 13:41:51.733191"
$ws.Range("F97").Value = "2021-10-05 13:41:51.733194"
$ws.Range("F98").Value = "2021-10-05 13:41:51.733197"
$ws.Range("F99").Value = "2021-10-05 13:41:51.733200"
